# October-2014-bank_statement.xlsx — "BalaRaju - Working on Spec"
#
# Replace row 2's record (Sekhar Beri) with Balaraju vankala's data, and
# drop row 3 (Priyanka Muddana) entirely so the sheet ends at row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: account number is cleared (blank), name + netpay updated; month
# ("October") in D2 is unchanged.
$ws.Range("A2").ClearContents()
$ws.Range("B2").Value = "Balaraju vankala"
$ws.Range("C2").Value = 87004.6

# Remove the old row 3 record completely (shifts nothing else up since it's
# the last row) so the used range becomes A1:D2.
$ws.Rows("3").Delete()

# Column A was sized (bestFit) to the old 16-digit account number; with that
# value gone the column narrows accordingly.
$ws.Columns("A").ColumnWidth = 13.375601926163727
